$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.103.60'
$ws.Range("E2").Value = '  +0.19%  '
$ws.Range("D3").Value = '1.838.65'
$ws.Range("E3").Value = '  +0.70%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '243.42'
$ws.Range("E5").Value = '  +0.86%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6294'
$ws.Range("E6").Value = '  -1.07%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.002'
$ws.Range("E7").Value = '  +0.13%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07577'
$ws.Range("E8").Value = '  +3.23%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2941'
$ws.Range("E9").Value = '  +0.20%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '22.63'
$ws.Range("E10").Value = '  -0.68%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07758'
$ws.Range("E11").Value = '  +1.17%  '
$ws.Range("D12").Value = '1.844.45'
$ws.Range("E12").Value = '  +0.98%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.975'
$ws.Range("E13").Value = '  -0.22%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6678'
$ws.Range("E14").Value = '  +0.73%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.00001005'
$ws.Range("E15").Value = '  +15.57%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '83.22'
$ws.Range("E16").Value = '  +1.50%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.087'
$ws.Range("E17").Value = '  +0.81%  '
$ws.Range("D18").Value = '29.112.01'
$ws.Range("E18").Value = '  +0.19%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '226.98'
$ws.Range("E19").Value = '  +0.79%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.45'
$ws.Range("E20").Value = '  +0.51%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.002'
$ws.Range("E21").Value = '  +0.17%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.228'
$ws.Range("E22").Value = '  +1.42%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.002'
$ws.Range("E23").Value = '  +0.11%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '160.02'
$ws.Range("E24").Value = '  +0.84%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1391'
$ws.Range("E25").Value = '  +1.73%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.518'
$ws.Range("E26").Value = '  +0.55%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.96'
$ws.Range("E27").Value = '  +0.38%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.500'
$ws.Range("E28").Value = '  -0.19%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.107'
$ws.Range("E29").Value = '  +0.42%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.024'
$ws.Range("E30").Value = '  -0.11%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.197'
$ws.Range("E31").Value = '  -0.50%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.05275'
$ws.Range("E32").Value = '  -0.66%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.851'
$ws.Range("E33").Value = '  +0.87%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7387'
$ws.Range("E34").Value = '  +0.20%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.139'
$ws.Range("E35").Value = '  -1.46%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.684'
$ws.Range("E36").Value = '  +1.24%  '
$ws.Range("D37").Value = '1.246.43'
$ws.Range("E37").Value = '  -3.91%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.764'
$ws.Range("E38").Value = '  +0.66%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01789'
$ws.Range("E39").Value = '  +0.09%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.381'
$ws.Range("E40").Value = '  +1.38%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9043'
$ws.Range("E41").Value = '  +0.61%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.003'
$ws.Range("E42").Value = '  +0.27%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '102.23'
$ws.Range("E43").Value = '  -0.39%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.00000000127'
$ws.Range("E44").Value = '  +6.20%  '
$ws.Range("D45").Value = '1.988.07'
$ws.Range("E45").Value = '  +0.71%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '64.45'
$ws.Range("E46").Value = '  +0.68%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5125'
$ws.Range("E47").Value = '  -0.23%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.4050'
$ws.Range("E48").Value = '  +1.62%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.952'
$ws.Range("E49").Value = '  +2.75%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05775'
$ws.Range("E50").Value = '  -0.53%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.726'
$ws.Range("E51").Value = '  +0.49%  '
